$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''298.16'
$ws.Range("E2").Value = '''-0.01%'
$ws.Range("D3").Value = '''31.29'
$ws.Range("E3").Value = '''-0.25%'
$ws.Range("D4").Value = '''5.096'
$ws.Range("E4").Value = '''-0.50%'
$ws.Range("D5").Value = '''0.08050'
$ws.Range("E5").Value = '''9.85%'
$ws.Range("D6").Value = '''2.529'
$ws.Range("E6").Value = '''44.61%'
$ws.Range("D7").Value = '''7.818'
$ws.Range("E7").Value = '''0.86%'
$ws.Range("D8").Value = '''3.812'
$ws.Range("E8").Value = '''2.29%'
$ws.Range("D9").Value = '''0.9179'
$ws.Range("E9").Value = '''-0.66%'
$ws.Range("D10").Value = '''0.1729'
$ws.Range("E10").Value = '''3.84%'
$ws.Range("E11").Value = '''6.16%'
$ws.Range("D12").Value = '''0.08653'
$ws.Range("E12").Value = '''8.97%'
$ws.Range("D13").Value = '''0.03029'
$ws.Range("E13").Value = '''1.32%'
$ws.Range("D14").Value = '''0.09965'
$ws.Range("E14").Value = '''0.53%'
$ws.Range("D15").Value = '''0.001490'
$ws.Range("E15").Value = '''-0.55%'
$ws.Range("D16").Value = '''0.005978'
$ws.Range("E16").Value = '''-4.38%'
$ws.Range("D17").Value = '''3.496'
$ws.Range("D18").Value = '''2.247'
$ws.Range("E18").Value = '''0.93%'
$ws.Range("E19").Value = '''1.80%'
$ws.Range("D20").Value = '''0.1337'
$ws.Range("E20").Value = '''1.62%'
$ws.Range("D21").Value = '''4.588'
$ws.Range("E21").Value = '''0.84%'
$ws.Range("D22").Value = '''0.1617'
$ws.Range("E22").Value = '''2.27%'
$ws.Range("D23").Value = '''0.04609'
$ws.Range("E23").Value = '''-0.89%'
$ws.Range("D24").Value = '''0.001249'
$ws.Range("E24").Value = '''2.95%'
$ws.Range("D25").Value = '''0.004434'
$ws.Range("E25").Value = '''-6.45%'
$ws.Range("D26").Value = '''0.0001201'
$ws.Range("E26").Value = '''-7.41%'
$ws.Range("D27").Value = '''0.0003429'
$ws.Range("E27").Value = '''83.26%'
$ws.Range("D39").Value = '''0.01793'
$ws.Range("E39").Value = '''3.70%'
$ws.Range("D40").Value = '''0.04524'
$ws.Range("E40").Value = '''1.58%'
$ws.Range("D41").Value = '''0.007030'
$ws.Range("E41").Value = '''-1.37%'
$ws.Range("D42").Value = '''0.1343'
$ws.Range("E42").Value = '''0.90%'
$ws.Range("D43").Value = '''0.002241'
$ws.Range("E43").Value = '''1.67%'
$ws.Range("D44").Value = '''0.009817'
$ws.Range("E44").Value = '''-9.08%'
$ws.Range("D45").Value = '''0.00006596'
$ws.Range("E45").Value = '''6.38%'
$ws.Range("E46").Value = '''0.04%'
$ws.Range("B47").Value = 'BOLO'
$ws.Range("C47").Value = 'https://coinranking.com/coin/ogrGe0dEab+bolo-bolo'
$ws.Range("D47").Value = '''0.8206'
$ws.Range("E47").Value = '''11.07%'
$ws.Range("B48").Value = 'CoinbaseStockToken'
$ws.Range("C48").Value = 'https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin'
$ws.Range("D48").Value = '''0.005238'
$ws.Range("E48").Value = '''-48.72%'
$ws.Range("D49").Value = '''0.00002101'
$ws.Range("E49").Value = '''0.04%'
$ws.Range("D50").Value = '''0.0002001'
$ws.Range("E50").Value = '''0.12%'